$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save the original row 3 values (to be moved to row 4)
$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$e3 = $ws.Range("E3").Value2
$f3 = $ws.Range("F3").Value2
$g3 = $ws.Range("G3").Value2
$h3 = $ws.Range("H3").Value2
$q3 = $ws.Range("Q3").Value2
$r3 = $ws.Range("R3").Value2

# Save the original row 4 values (to be moved to row 3)
$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2
$e4 = $ws.Range("E4").Value2
$f4 = $ws.Range("F4").Value2
$g4 = $ws.Range("G4").Value2
$h4 = $ws.Range("H4").Value2
$q4 = $ws.Range("Q4").Value2
$r4 = $ws.Range("R4").Value2

# Write row 4's original values into row 3
$ws.Range("A3").Value2 = $a4
$ws.Range("B3").Value2 = $b4
$ws.Range("E3").Value2 = $e4
$ws.Range("F3").Value2 = $f4
$ws.Range("G3").Value2 = $g4
$ws.Range("H3").Value2 = $h4
$ws.Range("Q3").Value2 = $q4
$ws.Range("R3").Value2 = $r4

# Write row 3's original values into row 4
$ws.Range("A4").Value2 = $a3
$ws.Range("B4").Value2 = $b3
$ws.Range("E4").Value2 = $e3
$ws.Range("F4").Value2 = $f3
$ws.Range("G4").Value2 = $g3
$ws.Range("H4").Value2 = $h3
$ws.Range("Q4").Value2 = $q3
$ws.Range("R4").Value2 = $r3
